$wb = $excel.ActiveWorkbook

# Add the new "StoreLocation" sheet after the last existing sheet ("Language")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "StoreLocation"

# Populate data matching the "Find a store" style rows used on other sheets
$newSheet.Range("A2").Value = 11102
$newSheet.Range("B2").Value = "New York, NY"

# Match the selection state captured in the target worksheet
$newSheet.Range("B2").Select()
